$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 304.09525
$ws.Range("I5").Value = 239.875
$ws.Range("J5").Value = 509.6
$ws.Range("K5").Value = 239.875
$ws.Range("L5").Value = 509.6
$ws.Range("M5").Value = -124.875
$ws.Range("N5").Value = -739.6

$ws.Range("H17").Value = 1292.25
$ws.Range("J17").Value = 1366.6666
$ws.Range("L17").Value = 4099.9998
$ws.Range("N17").Value = -4435.9998

$ws.Range("H28").Value = 936.8570999999999
$ws.Range("I28").Value = 936.8570999999999
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 936.8570999999999
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = -451.8570999999999
$ws.Range("N28").Value = ""

$ws.Range("H41").Value = 797.7143
$ws.Range("I41").Value = 422.5
$ws.Range("J41").Value = 947.8
$ws.Range("K41").Value = 422.5
$ws.Range("L41").Value = 947.8
$ws.Range("M41").Value = 17.5
$ws.Range("N41").Value = -1827.8

$ws.Range("H88").Value = 5004
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 5004
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 5004
$ws.Range("M88").Value = ""
$ws.Range("N88").Value = -5816

$ws.Range("H91").Value = 5004
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 5004
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 5004
$ws.Range("M91").Value = ""
$ws.Range("N91").Value = -7812

$ws.Range("H116").Value = 7042.7144
$ws.Range("I116").Value = 6716.5
$ws.Range("J116").Value = 9000
$ws.Range("K116").Value = 6716.5
$ws.Range("L116").Value = 9000
$ws.Range("M116").Value = -3274.5
$ws.Range("N116").Value = -15884

$ws.Range("H125").Value = 31252924
$ws.Range("J125").Value = 5268
$ws.Range("L125").Value = 47412
$ws.Range("N125").Value = -52332

$ws.Range("H132").Value = 2837.4614
$ws.Range("J132").Value = 7500
$ws.Range("L132").Value = 22500
$ws.Range("N132").Value = -27560

$ws.Range("H135").Value = 870.5
$ws.Range("I135").Value = 870.5
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 7834.5
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -5299.5
$ws.Range("N135").Value = ""

$ws.Range("H138").Value = 10493.228
$ws.Range("J138").Value = 10939.417
$ws.Range("L138").Value = 32818.251
$ws.Range("N138").Value = -43098.251

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7656.5947
$ws.Range("I32").Value = 6665.6284
$ws.Range("K32").Value = 6665.6284
$ws.Range("M32").Value = -6378.6284

$ws.Range("H45").Value = 2310.6924
$ws.Range("I45").Value = 1766.375
$ws.Range("J45").Value = 3181.6
$ws.Range("K45").Value = 1766.375
$ws.Range("L45").Value = 3181.6
$ws.Range("M45").Value = -1389.375
$ws.Range("N45").Value = -3935.6

$ws.Range("H61").Value = 5999
$ws.Range("J61").Value = 7999
$ws.Range("L61").Value = 7999
$ws.Range("N61").Value = -8423

$ws.Range("H74").Value = 3074.25
$ws.Range("I74").Value = 3074.25
$ws.Range("K74").Value = 3074.25
$ws.Range("M74").Value = -2200.25

$ws.Range("H77").Value = 3074.25
$ws.Range("I77").Value = 3074.25
$ws.Range("K77").Value = 15371.25
$ws.Range("M77").Value = -11003.25

$ws.Range("H97").Value = 506
$ws.Range("I97").Value = 506
$ws.Range("K97").Value = 506
$ws.Range("M97").Value = -10

$ws.Range("H122").Value = 1655
$ws.Range("I122").Value = 1655
$ws.Range("K122").Value = 4965
$ws.Range("M122").Value = -2515

$ws.Range("H132").Value = 2980.7144
$ws.Range("I132").Value = 2227.5
$ws.Range("K132").Value = 6682.5
$ws.Range("M132").Value = -4152.5

$ws.Range("H136").Value = 5999
$ws.Range("J136").Value = 7999
$ws.Range("L136").Value = 23997
$ws.Range("N136").Value = -29097

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 11219
$ws.Range("J4").Value = 13033.333
$ws.Range("L4").Value = 13033.333
$ws.Range("N4").Value = -13257.333

$ws.Range("H25").Value = 1412.5
$ws.Range("I25").Value = 1250
$ws.Range("J25").Value = 1900
$ws.Range("K25").Value = 1250
$ws.Range("L25").Value = 1900
$ws.Range("M25").Value = -1076
$ws.Range("N25").Value = -2248

$ws.Range("H31").Value = 3204
$ws.Range("I31").Value = 3014
$ws.Range("J31").Value = 4249
$ws.Range("K31").Value = 3014
$ws.Range("L31").Value = 4249
$ws.Range("M31").Value = -2719
$ws.Range("N31").Value = -4839

$ws.Range("H34").Value = 3204
$ws.Range("I34").Value = 3014
$ws.Range("J34").Value = 4249
$ws.Range("K34").Value = 3014
$ws.Range("L34").Value = 4249
$ws.Range("M34").Value = -2812
$ws.Range("N34").Value = -4653

$ws.Range("H41").Value = 8413.200000000001
$ws.Range("I41").Value = 3016.5
$ws.Range("K41").Value = 3016.5
$ws.Range("M41").Value = -2588.5

$ws.Range("H58").Value = 924
$ws.Range("I58").Value = 924
$ws.Range("K58").Value = 924
$ws.Range("M58").Value = -721

$ws.Range("H86").Value = 23239650
$ws.Range("I86").Value = 34854484
$ws.Range("K86").Value = 34854484
$ws.Range("M86").Value = -34853361

$ws.Range("H89").Value = 23239650
$ws.Range("I89").Value = 34854484
$ws.Range("K89").Value = 174272420
$ws.Range("M89").Value = -174266804

$ws.Range("H94").Value = 4551.125
$ws.Range("I94").Value = 4599
$ws.Range("K94").Value = 4599
$ws.Range("M94").Value = -4148

$ws.Range("H105").Value = 1031.9524
$ws.Range("I105").Value = 1025.3158
$ws.Range("J105").Value = 1095
$ws.Range("K105").Value = 1025.3158
$ws.Range("L105").Value = 1095
$ws.Range("M105").Value = 721.6841999999999
$ws.Range("N105").Value = -4589

$ws.Range("H136").Value = 924
$ws.Range("I136").Value = 924
$ws.Range("K136").Value = 2772
$ws.Range("M136").Value = -222

$ws.Range("H140").Value = 0
$ws.Range("I140").Value = 0
$ws.Range("K140").Value = 0
$ws.Range("M140").Value = ""

$ws.Range("H141").Value = 464781.1
$ws.Range("J141").Value = 464781.1
$ws.Range("L141").Value = 464781.1
$ws.Range("N141").Value = -475141.1

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 3799.1667
$ws.Range("I3").Value = 3799.1667
$ws.Range("K3").Value = 11397.5001
$ws.Range("M3").Value = -11285.5001

$ws.Range("H98").Value = 664.25
$ws.Range("I98").Value = 667.75
$ws.Range("J98").Value = 660.75
$ws.Range("K98").Value = 2003.25
$ws.Range("L98").Value = 1982.25
$ws.Range("M98").Value = -505.25
$ws.Range("N98").Value = -4978.25

$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").Value = ""

$ws.Range("H129").Value = 4998.8
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 4998.8
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 14996.4
$ws.Range("M129").Value = ""
$ws.Range("N129").Value = -24996.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1677
$ws.Range("I97").Value = 855
$ws.Range("K97").Value = 855
$ws.Range("M97").Value = -359

$ws.Range("H117").Value = 64450
$ws.Range("J117").Value = 64450
$ws.Range("L117").Value = 64450
$ws.Range("N117").Value = -71334

$ws.Range("H122").Value = 3522.8096
$ws.Range("J122").Value = 2313
$ws.Range("L122").Value = 6939
$ws.Range("N122").Value = -11839

$ws.Range("H132").Value = 4753.8335
$ws.Range("I132").Value = 2631
$ws.Range("K132").Value = 7893
$ws.Range("M132").Value = -5363

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 4325.4443
$ws.Range("I136").Value = 3989.8572
$ws.Range("K136").Value = 11969.5716
$ws.Range("M136").Value = -9419.571599999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3229.8823
$ws.Range("I132").Value = 3290.6
$ws.Range("K132").Value = 9871.799999999999
$ws.Range("M132").Value = -7341.799999999999

$ws.Range("H136").Value = 28919.158
$ws.Range("J136").Value = 5757
$ws.Range("L136").Value = 17271
$ws.Range("N136").Value = -22371
